$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44188
$ws.Range("J2").Value = 210
$ws.Range("K2").Value = 5000
$ws.Range("M2").Value = 5500
$ws.Range("O2").Value = "Provincia de Quillota"
$ws.Range("P2").Value = 344

# Row 3
$ws.Range("D3").Value = 44186
$ws.Range("J3").Value = 160

# Row 5
$ws.Range("D5").Value = 44230
$ws.Range("J5").Value = 250

# Row 6
$ws.Range("D6").Value = 44189
$ws.Range("J6").Value = 250

# Row 8
$ws.Range("D8").Value = 44232
$ws.Range("J8").Value = 250

# Row 9
$ws.Range("D9").Value = 44187
$ws.Range("J9").Value = 160

# Row 10
$ws.Range("D10").Value = 44204
$ws.Range("J10").Value = 430

# Row 11
$ws.Range("D11").Value = 44210
$ws.Range("J11").Value = 340

# Row 12
$ws.Range("D12").Value = 44292
$ws.Range("J12").Value = 90
$ws.Range("K12").Value = 6000
$ws.Range("M12").Value = 6000
$ws.Range("O12").Value = "Región Metropolitana"
$ws.Range("P12").Value = 375

# Row 13
$ws.Range("D13").Value = 44231
$ws.Range("J13").Value = 250

# Row 14
$ws.Range("D14").Value = 44215
